$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Correct" header in F1
$ws.Range("F1").Value = "Correct"

# Row 3 and Row 11 had their A:E contents shifted left by one column
# (question text moved from B->A position), so re-write them fully.
$ws.Range("A3").Value = "Which gland releases the Growth Hormone?"
$ws.Range("B3").Value = "Pineal gland"
$ws.Range("C3").Value = "Pituatary gland"
$ws.Range("D3").Value = "Adrenal glands"
$ws.Range("E3").Value = "Pancreas"

$ws.Range("A11").Value = "What processes occur in the Mitochondria of a cell?"
$ws.Range("B11").Value = "Urea cycle"
$ws.Range("C11").Value = "Glycolysis"
$ws.Range("D11").Value = "Boiling"
$ws.Range("E11").Value = "Haematopoiesis"

# D10 was the text "π" - fix it to be a plain number
$ws.Range("D10").Value = 2.43562235

# Fill the new "Correct" answer column (F2:F11) with the correct option letters
$ws.Range("F2").Value = "C"
$ws.Range("F3").Value = "B"
$ws.Range("F4").Value = "B"
$ws.Range("F5").Value = "D"
$ws.Range("F6").Value = "A"
$ws.Range("F7").Value = "B"
$ws.Range("F8").Value = "D"
$ws.Range("F9").Value = "A"
$ws.Range("F10").Value = "C"
$ws.Range("F11").Value = "A"

# Update the active selection to match the new state
$ws.Range("B16").Select()
